# psa_water_sensor_install / psa_water_sensor_uninstall
# Updated wording around manually typing of barcodes for nodes and gateways.
#
# The "manual entry" note that used to read:
#   "IMPORTANT: Double check correct numbers"
# is reworded to:
#   "Manually type here. IMPORTANT: Double check correct numbers"
# and is now used consistently for every "alt_*" manual-entry row (gateway,
# cover-crop node rep1/rep2, bare node rep1/rep2) - including the bare-node
# rep1 row, which previously carried its own one-off long-winded note
# ("Try again. If it fails multiple times you can manually type the barcode
# here. IMPORTANT: Double check and use correct cases and dashes.") that is
# now retired in favour of the shared, shorter note.

$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")

$newNote = "Manually type here. IMPORTANT: Double check correct numbers"

# Column C holds the "hint"/note label for each manual-entry row.
$survey.Cells.Item(15, 3).Value = $newNote   # alt_barcode_gateway hint
$survey.Cells.Item(22, 3).Value = $newNote   # alt_cover_crop_node_rep1 hint
$survey.Cells.Item(26, 3).Value = $newNote   # alt_bare_node_rep1 hint (was the long-winded one-off note)
$survey.Cells.Item(32, 3).Value = $newNote   # alt_cover_crop_node_rep2 hint
$survey.Cells.Item(36, 3).Value = $newNote   # alt_bare_node_rep2 hint

# Restore the view state left behind by the editor: back at the top of the
# sheet with E6 selected, rather than scrolled down with C39 selected.
$survey.Activate()
$survey.Range("E6").Select()
